# Aggiorno file need_to_buy.xlsx da R
# Refresh the rolling forecast table (A2:F15) with the latest values,
# shifting the date window forward and updating forecast columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45866
$ws.Range("B2").Value = 4654.8038038703
$ws.Range("C2").Value = 4578.42516404833
$ws.Range("D2").Value = 2376
$ws.Range("E2").Value = 5174.503849
$ws.Range("F2").Value = 113.421883715751

$ws.Range("A3").Value = 45867
$ws.Range("B3").Value = 4654.80380387156
$ws.Range("C3").Value = 4555.16558079577
$ws.Range("D3").Value = 2376
$ws.Range("E3").Value = 5174.503849
$ws.Range("F3").Value = 112.452734413509

$ws.Range("A4").Value = 45868
$ws.Range("B4").Value = 4654.80380387362
$ws.Range("C4").Value = 4577.28980553099
$ws.Range("D4").Value = 2376
$ws.Range("E4").Value = 5174.503849
$ws.Range("F4").Value = 113.374577110724

$ws.Range("A5").Value = 45869
$ws.Range("B5").Value = 4654.8038038703
$ws.Range("C5").Value = 4544.32462611719
$ws.Range("D5").Value = 2376
$ws.Range("E5").Value = 5174.503849
$ws.Range("F5").Value = 112.001027968621

$ws.Range("A6").Value = 45870
$ws.Range("B6").Value = 5180.74451596793
$ws.Range("C6").Value = 4219.81470893739
$ws.Range("D6").Value = 1944
$ws.Range("E6").Value = 5482.543494
$ws.Range("F6").Value = 107.400570290394

$ws.Range("A7").Value = 45871
$ws.Range("B7").Value = 908.164063516383
$ws.Range("C7").Value = 1557.79297321789
$ws.Range("D7").Value = 1944
$ws.Range("E7").Value = 1707.069026
$ws.Range("F7").Value = 17.1957473208962

$ws.Range("A8").Value = 45872
$ws.Range("B8").Value = 794.873906828036
$ws.Range("C8").Value = 1565.6422507254
$ws.Range("D8").Value = 1944
$ws.Range("E8").Value = 1582.629021
$ws.Range("F8").Value = 17.05822353739

$ws.Range("A9").Value = 45873
$ws.Range("B9").Value = 3947.0732721826
$ws.Range("C9").Value = 3909.21299633947
$ws.Range("D9").Value = 1944
$ws.Range("E9").Value = 4361.548203
$ws.Range("F9").Value = 99.1536636315361

$ws.Range("A10").Value = 45874
$ws.Range("B10").Value = 3947.0732721826
$ws.Range("C10").Value = 3920.13771345256
$ws.Range("D10").Value = 1944
$ws.Range("E10").Value = 4361.548203
$ws.Range("F10").Value = 99.6088601779148

$ws.Range("A11").Value = 45875
$ws.Range("B11").Value = 3947.0732721826
$ws.Range("C11").Value = 3920.98893648981
$ws.Range("D11").Value = 1944
$ws.Range("E11").Value = 4361.548203
$ws.Range("F11").Value = 99.644327804467

$ws.Range("A12").Value = 45876
$ws.Range("B12").Value = 3947.0732721826
$ws.Range("C12").Value = 3839.37853158575
$ws.Range("D12").Value = 1944
$ws.Range("E12").Value = 4361.548203
$ws.Range("F12").Value = 96.2438942667977

$ws.Range("A13").Value = 45877
$ws.Range("B13").Value = 3947.0732721826
$ws.Range("C13").Value = 3783.41647061902
$ws.Range("D13").Value = 1944
$ws.Range("E13").Value = 4361.548203
$ws.Range("F13").Value = 93.9121417265174

$ws.Range("A14").Value = 45878
$ws.Range("B14").Value = 719.129912541875
$ws.Range("C14").Value = 1405.81833555984
$ws.Range("D14").Value = 1944
$ws.Range("E14").Value = 1499.572444
$ws.Range("F14").Value = 10.0942027924151

$ws.Range("A15").Value = 45879
$ws.Range("B15").Value = 636.544690493214
$ws.Range("C15").Value = 1354.38935268938
$ws.Range("D15").Value = 1944
$ws.Range("E15").Value = 1408.859302
$ws.Range("F15").Value = 7.61266517484006

